{"js": "// 1) Drop the stray closing parenthesis: \"...External Customer.)\" -> \"...External Customer.\"\nconst target = \"as per request by the External Customer.)\";\nconst replacement = \"as per request by the External Customer.\";\n\nconst results = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the leftover \"_GoBack\" bookmark (the empty paragraph that hosted it stays put)\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmarkRange.load(\"isNullObject\");\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Drop the stray closing parenthesis: \"...External Customer.)\" -> \"...External Customer.\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"as per request by the External Customer.)\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"as per request by the External Customer.\"\n$rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, $null, $false, $null, 2) | Out-Null\n\n# 2) Remove the leftover \"_GoBack\" bookmark (the empty paragraph that hosted it stays put)\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
